$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logs")

# Insert a new column before column D (shifts old D..H to E..I)
$ws.Columns("D").Insert()

# --- Header row changes ---
$ws.Range("D1").Value = "sample(s)"

# --- Row 36-39 content updates (SingleR visualization rows) ---
$ws.Range("B36").Value = "as desc"
$ws.Range("E36").Value = "SCTv2 preselection aggrFalse meanLabelScore"
$ws.Range("B37").Value = "as desc"
$ws.Range("E37").Value = "SCTv2 preselection aggrTrue meanLabelScore"
$ws.Range("B38").Value = "as desc"
$ws.Range("E38").Value = "SCTv2 preselection aggrFalse maxLabelScore"
$ws.Range("B39").ClearContents()
$ws.Range("E39").Value = "SCTv2 preselection aggrTrue maxLabelScore"

# --- New rows 41-52 ---
$ws.Range("C41").Value = 'Notes'
$ws.Range("E41").Value = 'postselection data failing in visualization, probably wrong data in Kriegstein to SingleR so need to rerun (after organising parameters for visualization)'

$ws.Range("A42").Value = 'results'
$ws.Range("B42").Value = 'Pipe_SCTv2_23-06'
$ws.Range("C42").Value = 'individual'
$ws.Range("D42").Value = 'A'
$ws.Range("E42").Value = 'rerun'

$ws.Range("A43").Value = 'results'
$ws.Range("B43").Value = 'Pipe_SCTv2_23-06'
$ws.Range("C43").Value = 'individual'
$ws.Range("D43").Value = 'C'
$ws.Range("E43").Value = 'rerun'

$ws.Range("A44").Value = 'results'
$ws.Range("B44").Value = 'Pipe_SCTv2_23-06'
$ws.Range("C44").Value = 'individual'
$ws.Range("D44").Value = 'N'
$ws.Range("E44").Value = 'rerun'

$ws.Range("A45").Value = 'results'
$ws.Range("B45").Value = 'Pipe_SCTv2_23-06'
$ws.Range("C45").Value = 'integration'
$ws.Range("D45").Value = 'A + C'
$ws.Range("E45").Value = 'old selection'

$ws.Range("A46").Value = 'results'
$ws.Range("B46").Value = 'Pipe_SCTv2_23-06'
$ws.Range("C46").Value = 'integration'
$ws.Range("D46").Value = 'A + C'
$ws.Range("E46").Value = 'new selection'

$ws.Range("A47").Value = 'results'
$ws.Range("B47").Value = 'Pipe_SCTv2_23-06'
$ws.Range("C47").Value = 'integration'
$ws.Range("D47").Value = 'N + C'
$ws.Range("E47").Value = 'old selection'

$ws.Range("A48").Value = 'results'
$ws.Range("B48").Value = 'Pipe_SCTv2_23-06'
$ws.Range("C48").Value = 'integration'
$ws.Range("D48").Value = 'N + C'
$ws.Range("E48").Value = 'new selection'

$ws.Range("A49").Value = 'results'
$ws.Range("B49").Value = '2022-06-23 16-03-44'
$ws.Range("C49").Value = 'DEG'
$ws.Range("D49").Value = 'A + C'
$ws.Range("E49").Value = 'old selection'

$ws.Range("A50").Value = 'results'
$ws.Range("B50").Value = '2022-06-23 16-04-56'
$ws.Range("C50").Value = 'DEG'
$ws.Range("D50").Value = 'A + C'
$ws.Range("E50").Value = 'new selection'

$ws.Range("C51").Value = 'DEG'
$ws.Range("D51").Value = 'N + C'
$ws.Range("E51").Value = 'old selection'

$ws.Range("C52").Value = 'DEG'
$ws.Range("D52").Value = 'N + C'
$ws.Range("E52").Value = 'new selection'

# --- View state updates on "logs" sheet ---
$excel.ActiveWindow.ScrollRow = 28
$ws.Range("B50").Select()

# --- View state updates on "SCTv2 pipeline runtime" sheet ---
$ws2 = $wb.Worksheets.Item("SCTv2 pipeline runtime")
$ws2.Activate()
$ws2.Range("C20").Select()
$ws.Activate()
